# Weekly data refresh: insert a new daily price record for
# "Vega Modelo de Temuco - Apio" at row 244, pushing all the existing
# records (old rows 244-312) down by one row (new rows 245-313).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 244; Excel shifts rows 244:312 down
# to 245:313 (same behaviour as right-clicking the row header -> Insert).
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with the new record. The
# categorical columns (market/region/product/etc.) mirror the rest of
# the sheet's "Provincia del Elquí" / "Primera" entries; only the
# date, volume and price columns carry new values.
$ws.Cells.Item(244, 1).Value  = 10
$ws.Cells.Item(244, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(244, 3).Value  = "La Araucanía"
$ws.Cells.Item(244, 4).Value  = 44722
$ws.Cells.Item(244, 5).Value  = 9
$ws.Cells.Item(244, 6).Value  = 100112017
$ws.Cells.Item(244, 7).Value  = "Apio"
$ws.Cells.Item(244, 8).Value  = "Americana (o)"
$ws.Cells.Item(244, 9).Value  = "Primera"
$ws.Cells.Item(244, 10).Value = 50
$ws.Cells.Item(244, 11).Value = 9000
$ws.Cells.Item(244, 12).Value = 9000
$ws.Cells.Item(244, 13).Value = 9000
$ws.Cells.Item(244, 14).Value = "$/docena de matas"
$ws.Cells.Item(244, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(244, 16).Value = 1500
$ws.Cells.Item(244, 17).Value = 6
$ws.Cells.Item(244, 18).Value = "Hortaliza"
